$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 549.5925999999999
$ws.Range("J17").Value = 566.88464
$ws.Range("L17").Value = 1700.65392
$ws.Range("N17").Value = -2036.65392
$ws.Range("H40").Value = 10006249
$ws.Range("J40").Value = 33346666
$ws.Range("L40").Value = 33346666
$ws.Range("N40").Value = -33347016
$ws.Range("H76").Value = 4654.25
$ws.Range("I76").Value = 4512.3335
$ws.Range("K76").Value = 4512.3335
$ws.Range("M76").Value = -4197.3335
$ws.Range("H79").Value = 4654.25
$ws.Range("I79").Value = 4512.3335
$ws.Range("K79").Value = 4512.3335
$ws.Range("M79").Value = -3420.3335
$ws.Range("H88").Value = 914281.4
$ws.Range("J88").Value = 1255231.6
$ws.Range("L88").Value = 1255231.6
$ws.Range("N88").Value = -1256043.6
$ws.Range("H91").Value = 914281.4
$ws.Range("J91").Value = 1255231.6
$ws.Range("L91").Value = 1255231.6
$ws.Range("N91").Value = -1258039.6
$ws.Range("H106").Value = 12246.2
$ws.Range("I106").Value = 13559.064
$ws.Range("K106").Value = 13559.064
$ws.Range("M106").Value = -12928.064
$ws.Range("H116").Value = 4391.0435
$ws.Range("I116").Value = 4443.8
$ws.Range("J116").Value = 4350.4614
$ws.Range("K116").Value = 4443.8
$ws.Range("L116").Value = 4350.4614
$ws.Range("M116").Value = -1001.8
$ws.Range("N116").Value = -11234.4614
$ws.Range("H132").Value = 3127.2222
$ws.Range("I132").Value = 2828.7058
$ws.Range("K132").Value = 8486.117400000001
$ws.Range("M132").Value = -5956.117400000001
$ws.Range("H135").Value = 1872.3182
$ws.Range("I135").Value = 1872.3182
$ws.Range("K135").Value = 16850.8638
$ws.Range("M135").Value = -14315.8638
$ws.Range("H137").Value = 1041.6666
$ws.Range("I137").Value = 920.3333
$ws.Range("K137").Value = 2760.9999
$ws.Range("M137").Value = -210.9998999999998
$ws.Range("H138").Value = 3949.1853
$ws.Range("I138").Value = 2523
$ws.Range("J138").Value = 5090.1333
$ws.Range("K138").Value = 7569
$ws.Range("L138").Value = 15270.3999
$ws.Range("M138").Value = -2429
$ws.Range("N138").Value = -25550.3999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2368.35
$ws.Range("I2").Value = 1506.8
$ws.Range("K2").Value = 1506.8
$ws.Range("M2").Value = -1393.8
$ws.Range("H32").Value = 1860.0339
$ws.Range("I32").Value = 886.34546
$ws.Range("K32").Value = 886.34546
$ws.Range("M32").Value = -599.34546
$ws.Range("H45").Value = 1984.6666
$ws.Range("I45").Value = 2302.5
$ws.Range("J45").Value = 1349
$ws.Range("K45").Value = 2302.5
$ws.Range("L45").Value = 1349
$ws.Range("M45").Value = -1925.5
$ws.Range("N45").Value = -2103
$ws.Range("H61").Value = 3089.1428
$ws.Range("I61").Value = 2723.3125
$ws.Range("J61").Value = 4259.8
$ws.Range("K61").Value = 2723.3125
$ws.Range("L61").Value = 4259.8
$ws.Range("M61").Value = -2511.3125
$ws.Range("N61").Value = -4683.8
$ws.Range("H102").Value = 5727.4
$ws.Range("I102").Value = 5408
$ws.Range("J102").Value = 7005
$ws.Range("K102").Value = 5408
$ws.Range("L102").Value = 7005
$ws.Range("M102").Value = -3786
$ws.Range("N102").Value = -10249
$ws.Range("H116").Value = 2368.35
$ws.Range("I116").Value = 1506.8
$ws.Range("K116").Value = 1506.8
$ws.Range("M116").Value = 787.2
$ws.Range("H122").Value = 5004.6665
$ws.Range("J122").Value = 4507.25
$ws.Range("L122").Value = 13521.75
$ws.Range("N122").Value = -18421.75
$ws.Range("H125").Value = 79999
$ws.Range("J125").Value = 79999
$ws.Range("L125").Value = 79999
$ws.Range("N125").Value = -89839
$ws.Range("H136").Value = 3089.1428
$ws.Range("I136").Value = 2723.3125
$ws.Range("J136").Value = 4259.8
$ws.Range("K136").Value = 8169.9375
$ws.Range("L136").Value = 12779.4
$ws.Range("M136").Value = -5619.9375
$ws.Range("N136").Value = -17879.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2368.35
$ws.Range("I3").Value = 1506.8
$ws.Range("K3").Value = 1506.8
$ws.Range("M3").Value = -1392.8
$ws.Range("H86").Value = 2836.5
$ws.Range("J86").Value = 3032
$ws.Range("L86").Value = 3032
$ws.Range("N86").Value = -5278
$ws.Range("H89").Value = 2836.5
$ws.Range("J89").Value = 3032
$ws.Range("L89").Value = 15160
$ws.Range("N89").Value = -26392
$ws.Range("H99").Value = 2173.0715
$ws.Range("I99").Value = 2186.3845
$ws.Range("K99").Value = 2186.3845
$ws.Range("M99").Value = -688.3845000000001
$ws.Range("H107").Value = 1535.4667
$ws.Range("I107").Value = 1241
$ws.Range("K107").Value = 1241
$ws.Range("M107").Value = 679
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2639.7334
$ws.Range("I16").Value = 5349.2
$ws.Range("K16").Value = 5349.2
$ws.Range("M16").Value = -5062.2
$ws.Range("H31").Value = 6717.15
$ws.Range("J31").Value = 6614.4116
$ws.Range("L31").Value = 6614.4116
$ws.Range("N31").Value = -7204.4116
$ws.Range("H34").Value = 6717.15
$ws.Range("J34").Value = 6614.4116
$ws.Range("L34").Value = 6614.4116
$ws.Range("N34").Value = -7018.4116
$ws.Range("H113").Value = 2639.7334
$ws.Range("I113").Value = 5349.2
$ws.Range("K113").Value = 5349.2
$ws.Range("M113").Value = -3179.2
$ws.Range("H132").Value = 2150.3684
$ws.Range("I132").Value = 2085.7646
$ws.Range("J132").Value = 2699.5
$ws.Range("K132").Value = 6257.293799999999
$ws.Range("L132").Value = 8098.5
$ws.Range("M132").Value = -3727.293799999999
$ws.Range("N132").Value = -13158.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 869666.4
$ws.Range("I2").Value = 1111146
$ws.Range("J2").Value = 339.8
$ws.Range("K2").Value = 6666876
$ws.Range("L2").Value = 2038.8
$ws.Range("M2").Value = -6666763
$ws.Range("N2").Value = -2264.8
$ws.Range("H68").Value = 2148.6667
$ws.Range("J68").Value = 2651.1428
$ws.Range("L68").Value = 7953.428400000001
$ws.Range("N68").Value = -9575.428400000001
$ws.Range("H69").Value = 4816
$ws.Range("I69").Value = 4816
$ws.Range("K69").Value = 14448
$ws.Range("M69").Value = -13637
$ws.Range("H71").Value = 2148.6667
$ws.Range("J71").Value = 2651.1428
$ws.Range("L71").Value = 23860.2852
$ws.Range("N71").Value = -31972.2852
$ws.Range("H72").Value = 4816
$ws.Range("I72").Value = 4816
$ws.Range("K72").Value = 43344
$ws.Range("M72").Value = -39288
$ws.Range("H107").Value = 1205.8379
$ws.Range("I107").Value = 703.7143
$ws.Range("J107").Value = 1864.875
$ws.Range("K107").Value = 2111.1429
$ws.Range("L107").Value = 5594.625
$ws.Range("M107").Value = -191.1428999999998
$ws.Range("N107").Value = -9434.625
$ws.Range("H115").Value = 101186.43
$ws.Range("I115").Value = 117976
$ws.Range("J115").Value = 449
$ws.Range("K115").Value = 353928
$ws.Range("L115").Value = 1347
$ws.Range("M115").Value = -352753
$ws.Range("N115").Value = -3697
$ws.Range("H131").Value = 9093127
$ws.Range("I131").Value = 7693837.5
$ws.Range("J131").Value = 11114324
$ws.Range("K131").Value = 23081512.5
$ws.Range("L131").Value = 33342972
$ws.Range("M131").Value = -23076472.5
$ws.Range("N131").Value = -33353052
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1382.3158
$ws.Range("I107").Value = 1598.6154
$ws.Range("K107").Value = 1598.6154
$ws.Range("M107").Value = 321.3846000000001
$ws.Range("H113").Value = 3403.818
$ws.Range("I113").Value = 3170.2144
$ws.Range("K113").Value = 3170.2144
$ws.Range("M113").Value = -1000.2144
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 33338982
$ws.Range("I7").Value = 5603.1665
$ws.Range("J7").Value = 166672500
$ws.Range("K7").Value = 5603.1665
$ws.Range("L7").Value = 166672500
$ws.Range("M7").Value = -5491.1665
$ws.Range("N7").Value = -166672724
$ws.Range("H16").Value = 2821.25
$ws.Range("I16").Value = 2392
$ws.Range("J16").Value = 3035.875
$ws.Range("K16").Value = 2392
$ws.Range("L16").Value = 3035.875
$ws.Range("M16").Value = -2222
$ws.Range("N16").Value = -3375.875
$ws.Range("H122").Value = 6755.125
$ws.Range("I122").Value = 7208.6
$ws.Range("K122").Value = 21625.8
$ws.Range("M122").Value = -19175.8
$ws.Range("H126").Value = 33338982
$ws.Range("I126").Value = 5603.1665
$ws.Range("J126").Value = 166672500
$ws.Range("K126").Value = 16809.4995
$ws.Range("L126").Value = 500017500
$ws.Range("M126").Value = -14339.4995
$ws.Range("N126").Value = -500022440
$ws.Range("H132").Value = 7637.023
$ws.Range("I132").Value = 6908.8438
$ws.Range("J132").Value = 9578.833000000001
$ws.Range("K132").Value = 20726.5314
$ws.Range("L132").Value = 28736.499
$ws.Range("M132").Value = -18196.5314
$ws.Range("N132").Value = -33796.499
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5025.75
$ws.Range("I2").Value = 5025.75
$ws.Range("K2").Value = 5025.75
$ws.Range("M2").Value = -4913.75
$ws.Range("H126").Value = 8774781
$ws.Range("I126").Value = 10755349
$ws.Range("J126").Value = 3698.5715
$ws.Range("K126").Value = 32266047
$ws.Range("L126").Value = 11095.7145
$ws.Range("M126").Value = -32263577
$ws.Range("N126").Value = -16035.7145
$ws.Range("H132").Value = 3298.9092
$ws.Range("I132").Value = 3170.2856
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 9510.856800000001
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -6980.856800000001
$ws.Range("N132").Value = -23060
